# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted for "Agrícola del Norte S.A. de Arica - Maracuyá"
# This pushes all existing data rows (57-137) down by one (to 58-138) and
# populates the newly opened row 57 with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting rows 57:137 down to 58:138
$ws.Rows.Item(57).Insert()

# Populate the new row 57 with the new weekly record
$ws.Cells.Item(57, 1).Value2 = 1
$ws.Cells.Item(57, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value2 = 44771
$ws.Cells.Item(57, 5).Value2 = 15
$ws.Cells.Item(57, 6).Value2 = "Fruta"
$ws.Cells.Item(57, 7).Value2 = 100108
$ws.Cells.Item(57, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(57, 9).Value2 = 100108003
$ws.Cells.Item(57, 10).Value2 = "Maracuyá"
$ws.Cells.Item(57, 11).Value2 = "Sin especificar"
$ws.Cells.Item(57, 12).Value2 = "Primera"
$ws.Cells.Item(57, 13).Value2 = 130
$ws.Cells.Item(57, 14).Value2 = 24000
$ws.Cells.Item(57, 15).Value2 = 25000
$ws.Cells.Item(57, 16).Value2 = 24500
$ws.Cells.Item(57, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(57, 18).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 19).Value2 = 1225
$ws.Cells.Item(57, 20).Value2 = 20
